$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-18 Sunday", 2)

$d.Content.Find.Execute("72×33=2376", $true, $false, $false, $false, $false, $true, 1, $false, "83×53=4399", 2)
$d.Content.Find.Execute("50×60=3000", $true, $false, $false, $false, $false, $true, 1, $false, "93×90=8370", 2)
$d.Content.Find.Execute("17×62=1054", $true, $false, $false, $false, $false, $true, 1, $false, "56×24=1344", 2)
$d.Content.Find.Execute("49×70=3430", $true, $false, $false, $false, $false, $true, 1, $false, "18×62=1116", 2)
$d.Content.Find.Execute("47×48=2256", $true, $false, $false, $false, $false, $true, 1, $false, "97×57=5529", 2)

$d.Content.Find.Execute("85×42=3570", $true, $false, $false, $false, $false, $true, 1, $false, "55×74=4070", 2)
$d.Content.Find.Execute("12×74=888", $true, $false, $false, $false, $false, $true, 1, $false, "27×55=1485", 2)
$d.Content.Find.Execute("17×61=1037", $true, $false, $false, $false, $false, $true, 1, $false, "55×68=3740", 2)
$d.Content.Find.Execute("18×91=1638", $true, $false, $false, $false, $false, $true, 1, $false, "18×74=1332", 2)
$d.Content.Find.Execute("31×78=2418", $true, $false, $false, $false, $false, $true, 1, $false, "30×48=1440", 2)

$d.Content.Find.Execute("41×68=2788", $true, $false, $false, $false, $false, $true, 1, $false, "92×51=4692", 2)
$d.Content.Find.Execute("38×84=3192", $true, $false, $false, $false, $false, $true, 1, $false, "31×66=2046", 2)
$d.Content.Find.Execute("75×83=6225", $true, $false, $false, $false, $false, $true, 1, $false, "85×63=5355", 2)
$d.Content.Find.Execute("34×48=1632", $true, $false, $false, $false, $false, $true, 1, $false, "27×40=1080", 2)
$d.Content.Find.Execute("73×35=2555", $true, $false, $false, $false, $false, $true, 1, $false, "26×25=650", 2)

$d.Content.Find.Execute("97×53=5141", $true, $false, $false, $false, $false, $true, 1, $false, "56×68=3808", 2)
$d.Content.Find.Execute("71×44=3124", $true, $false, $false, $false, $false, $true, 1, $false, "61×97=5917", 2)
$d.Content.Find.Execute("16×74=1184", $true, $false, $false, $false, $false, $true, 1, $false, "24×58=1392", 2)
$d.Content.Find.Execute("56×65=3640", $true, $false, $false, $false, $false, $true, 1, $false, "47×89=4183", 2)
$d.Content.Find.Execute("11×78=858", $true, $false, $false, $false, $false, $true, 1, $false, "43×19=817", 2)

$d.Content.Find.Execute("80×79=6320", $true, $false, $false, $false, $false, $true, 1, $false, "59×79=4661", 2)
$d.Content.Find.Execute("89×63=5607", $true, $false, $false, $false, $false, $true, 1, $false, "15×57=855", 2)
$d.Content.Find.Execute("41×19=779", $true, $false, $false, $false, $false, $true, 1, $false, "12×40=480", 2)
$d.Content.Find.Execute("69×73=5037", $true, $false, $false, $false, $false, $true, 1, $false, "64×65=4160", 2)
$d.Content.Find.Execute("59×61=3599", $true, $false, $false, $false, $false, $true, 1, $false, "27×28=756", 2)
